# Auto-generated script to apply cryptos.xlsx data refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.483.45"
$ws.Range("E2").Value = "  +5.12%  "
$ws.Range("D3").Value = "3.506.66"
$ws.Range("E3").Value = "  +2.32%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "418.29"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "132.26"
$ws.Range("E6").Value = "  +2.85%  "
$ws.Range("E7").Value = "  +4.61%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.776"
$ws.Range("E9").Value = "  +6.63%  "
$ws.Range("E10").Value = "  +16.05%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "43.22"
$ws.Range("E11").Value = "  +1.09%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000264"
$ws.Range("E12").Value = "  +20.32%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "9.98"
$ws.Range("E13").Value = "  +8.13%  "
$ws.Range("D14").Value = "4.067.64"
$ws.Range("E14").Value = "  +2.55%  "
$ws.Range("E15").Value = "  +0.20%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.45"
$ws.Range("E16").Value = "  +0.11%  "
$ws.Range("D17").Value = "3.525.44"
$ws.Range("E17").Value = "  +2.97%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.84"
$ws.Range("E18").Value = "  +1.30%  "
$ws.Range("E19").Value = "  +2.40%  "
$ws.Range("D20").Value = "65.452.05"
$ws.Range("E20").Value = "  +5.14%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "451.32"
$ws.Range("E21").Value = "  -4.67%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "90.31"
$ws.Range("E22").Value = "  -1.36%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.24"
$ws.Range("E23").Value = "  -0.60%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.26"
$ws.Range("E24").Value = "  +1.01%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.39"
$ws.Range("E25").Value = "  +3.34%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.91"
$ws.Range("E26").Value = "  +2.04%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "34.06"
$ws.Range("E27").Value = "  +1.75%  "
$ws.Range("B28").Value = "Toncoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.82"
$ws.Range("E28").Value = "  +7.21%  "
$ws.Range("B29").Value = "Cosmos"
$ws.Range("C29").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "12.53"
$ws.Range("E29").Value = "  +5.51%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.43"
$ws.Range("E30").Value = "  -3.84%  "
$ws.Range("E31").Value = "  +5.62%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.162"
$ws.Range("E32").Value = "  -1.41%  "
$ws.Range("B33").Value = "Dai"
$ws.Range("C33").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.998"
$ws.Range("E33").Value = "  -0.17%  "
$ws.Range("B34").Value = "InjectiveProtocol"
$ws.Range("C34").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "39.53"
$ws.Range("E34").Value = "  -3.12%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "57.45"
$ws.Range("E35").Value = "  -0.87%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0506"
$ws.Range("E36").Value = "  +3.87%  "
$ws.Range("D37").Value = "0.0₃0740"
$ws.Range("E37").Value = "  +36.06%  "
$ws.Range("E38").Value = "  +9.44%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.998"
$ws.Range("E39").Value = "  -0.08%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.05"
$ws.Range("E40").Value = "  +1.07%  "
$ws.Range("B41").Value = "WEMIXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.76"
$ws.Range("E41").Value = "  +4.23%  "
$ws.Range("B42").Value = "NEARProtocol"
$ws.Range("C42").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.50"
$ws.Range("E42").Value = "  +3.78%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "145.78"
$ws.Range("E43").Value = "  +0.25%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.29"
$ws.Range("E44").Value = "  -1.00%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.311"
$ws.Range("E45").Value = "  -3.26%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.00"
$ws.Range("E46").Value = "  -3.13%  "
$ws.Range("E47").Value = "  -2.15%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "15.77"
$ws.Range("E48").Value = "  -3.13%  "
$ws.Range("E49").Value = "  +3.52%  "
$ws.Range("E50").Value = "  +10.71%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "21.60"
$ws.Range("E51").Value = "  -3.27%  "
